$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a WordprocessingML block fragment (one or more <w:p> elements)
# in the minimal flat-OPC envelope that Range.InsertXML() expects.
# ---------------------------------------------------------------------------
function New-FlatOpc([string]$innerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body>' + $innerXml + '</w:body>' +
            '</w:document>' +
          '</pkg:xmlData>' +
        '</pkg:part>' +
      '</pkg:package>'
}

# Finds the first paragraph (searching from the top of the document) whose
# text contains $searchText, then replaces the *entire paragraph* with the
# literal XML supplied in $paraXml (a single <w:p>...</w:p> element).
function Set-ParagraphXmlByText([string]$searchText, [string]$paraXml) {
    $doc = $word.ActiveDocument
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText)
    if (-not $found) {
        throw "Set-ParagraphXmlByText: text not found: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $pRange = $para.Range
    $pRange.InsertXML((New-FlatOpc $paraXml))
}

# Replaces the single paragraph that lives inside a table cell with the
# literal XML supplied in $paraXml.
function Set-CellParagraphXml($cell, [string]$paraXml) {
    $pRange = $cell.Range.Paragraphs(1).Range
    $pRange.InsertXML((New-FlatOpc $paraXml))
}

# ===========================================================================
# 1) Modificabilidad, bullet 1 (AT002) - body list item.
#    Split the sentence into two runs ("Facilidad para " + the new wording)
#    and relocate the _GoBack bookmark to the end of this paragraph.
# ===========================================================================
$at002Body = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
             '<w:r><w:t xml:space="preserve">Facilidad para </w:t></w:r>' +
             '<w:r><w:t>que el sistema use nuevos formatos de ficheros.</w:t></w:r>' +
             '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
             '</w:p>'
Set-ParagraphXmlByText "Facilidad para modificar, sustituir e intercambiar el sistema de parseado de ficheros, debido a la posibilidad de permitir varios formatos en un futuro." $at002Body

# ===========================================================================
# 2) Modificabilidad, bullet 2 (AT003) - body list item.
#    New wording as one run, plus a trailing-space run.
# ===========================================================================
$at003Body = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
             '<w:r><w:t>Facilidad para que el sistema use diferentes tipos de bases de datos</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '</w:p>'
Set-ParagraphXmlByText "Facilidad para modificar, sustituir e intercambiar el sistema de que lleva a cabo la comunicación con la base de datos, debido a la posibilidad de que la empresa decida usar un sistema diferente." $at003Body

# ===========================================================================
# 3) Testabilidad bullet (AT005) - body list item: reword, then drop the
#    blank paragraph that used to sit right after it.
# ===========================================================================
$at005Body = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
             '<w:r><w:t>Facilidad para probar el correcto procesado de las preguntas.</w:t></w:r>' +
             '</w:p>'
Set-ParagraphXmlByText "Facilidad para probar la fiabilidad del sistema, en especial en el sistema de conversión de preguntas." $at005Body

$d2 = $word.ActiveDocument
$rng2 = $d2.Content
$rng2.Find.Execute("Facilidad para probar el correcto procesado de las preguntas.") | Out-Null
$testabilidadPara = $rng2.Paragraphs(1)
$nextPara = $testabilidadPara.Next()
if ($nextPara.Range.Text.Trim() -eq "") {
    $nextPara.Range.Delete()
}

# ===========================================================================
# 4) Usabilidad bullet (AT006) - body list item: merge the two runs into one
#    with the new wording and drop the paragraph's <w:spacing w:after="0"/>.
# ===========================================================================
$at006Body = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
             '<w:r><w:t>Posibilidad de configurar los ficheros y formatos.</w:t></w:r>' +
             '</w:p>'
Set-ParagraphXmlByText "Facilidad para configurar y depurar el sistema por el usuario." $at006Body

# ===========================================================================
# 5) Remove the stale _GoBack bookmark left on the page-break paragraph (it
#    now lives solely inside the AT002 paragraph handled in step 1). That
#    paragraph is identified by containing a page-break character.
# ===========================================================================
$d3 = $word.ActiveDocument
for ($i = 1; $i -le $d3.Paragraphs.Count; $i++) {
    $para = $d3.Paragraphs($i)
    $r = $para.Range
    if ($r.Text.Contains([char]12)) {
        $r.InsertXML((New-FlatOpc '<w:p><w:r><w:br w:type="page"/></w:r></w:p>'))
        break
    }
}

# ===========================================================================
# 6) Table cell description text (column 2 = "Descripción").
# ===========================================================================
$d4 = $word.ActiveDocument
$tbl = $d4.Tables(1)

Set-CellParagraphXml $tbl.Cell(3, 2) '<w:p><w:r><w:t>Facilidad para que el sistema use nuevos formatos de ficheros.</w:t></w:r></w:p>'
Set-CellParagraphXml $tbl.Cell(4, 2) '<w:p><w:r><w:t>Facilidad para que el sistema use diferentes tipos de bases de datos</w:t></w:r></w:p>'
Set-CellParagraphXml $tbl.Cell(6, 2) '<w:p><w:r><w:t>Facilidad para probar el correcto procesado de las preguntas.</w:t></w:r></w:p>'
Set-CellParagraphXml $tbl.Cell(7, 2) '<w:p><w:r><w:t>Posibilidad de configurar los ficheros y formatos.</w:t></w:r></w:p>'
